# wwTHg_Comb_OutflowR.xlsx update
#
# The original workbook has a single sheet "4_wwTHg_Comb" with 31 days of
# wwTHg load data (rows 3-33). Row 21 (date serial 42721) carries an
# implausible Flow value (0.6, far below the surrounding ~200-13000 range)
# that was flagged during QA.
#
# This edit:
#   1. Renames the existing sheet to "4_wwTHg_Comb_31ct" (the "as received,
#      31 day count" version) and highlights the suspect row in light blue
#      so it's easy to spot.
#   2. Adds a second sheet "4_wwTHg_Comb_30ct", a duplicate of the first
#      with the suspect row removed entirely (the "cleaned, 30 day count"
#      version), tinted with the same light blue tab color, and leaves it
#      as the active/selected sheet - this is now the sheet analysts should
#      use going forward.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$highlight = 16444072   # RGB(168, 234, 250) == FFA8EAFA

# --- Build the cleaned "30ct" sheet as a copy of the original -----------
[void]$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# Drop the anomalous row (date 42721, Flow = 0.6) - rows below shift up.
$ws2.Rows.Item(21).Delete()

# --- Names ----------------------------------------------------------------
$ws1.Name = "4_wwTHg_Comb_31ct"
$ws2.Name = "4_wwTHg_Comb_30ct"

# --- Flag the suspect row on the original (31ct) sheet ---------------------
$ws1.Range("A21:D21").Interior.Color = $highlight

# --- Make the cleaned (30ct) sheet the active tab, tinted to match --------
$ws2.Tab.Color = $highlight
[void]$ws2.Range("E32").Select()
[void]$ws2.Activate()
